$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the indicator description text (B4) ---
$ws.Range("B4").Value = "4.1.1 Доля детей и молодежи, приходящаяся на a) учащихся 2 и 3 классов; b) выпускников начальной школы; и c) выпускников неполной средней школы, которые достигли по меньшей мере минимального уровня владения навыками i) чтения и i) математики, в разбивке по полу"
$ws.Rows.Item(4).RowHeight = 54
$ws.Range("B4").Font.Name = "Calibri"

# --- 2. Update organization contact person (B7) ---
$ws.Range("B7").Value = "1) Калымбетова Ы.И.`n2) Шамшидинова Бактыгуль Сабыржановна"
$ws.Range("B7").Font.Name = "Calibri"

# --- 3. Update contact email (B8) ---
$ws.Range("B8").Value = "1) yryskan.kalymbetova@gmail.com`n2) bakula68@mail.ru"
$ws.Range("B8").Font.Name = "Calibri"

# --- 4. Update contact phone (B9) ---
$ws.Range("B9").Value = "1) +996312 32 46 55`n2) +996312622962"
$ws.Range("B9").Font.Name = "Calibri"

# --- 5. Update organization website (B10) ---
$ws.Range("B10").Value = "1) www.stat.gov.kg`n2) http://ntc.kg/;  `nhttp://testing.kg"
$ws.Range("B10").Font.Name = "Calibri"

# --- 6. Move the active cell selection to B9 ---
$ws.Range("B9").Select() | Out-Null
